$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

$ws.Range("H2").Value = "S101"
$ws.Range("H3").Value = "S102"
$ws.Range("H4").Value = "S103"
$ws.Range("H5").Value = "S101"
$ws.Range("H6").Value = "S104"
$ws.Range("H7").Value = "S105"
$ws.Range("H8").Value = "S106"
$ws.Range("H9").Value = "S107"
$ws.Range("H10").Value = "S108"
$ws.Range("H11").Value = "S109"
